$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 807.9
$ws.Range("I2").Value = 234.58333
$ws.Range("J2").Value = 1190.1111
$ws.Range("K2").Value = 234.58333
$ws.Range("L2").Value = 1190.1111
$ws.Range("M2").Value = -121.58333
$ws.Range("N2").Value = -1416.1111
$ws.Range("H5").Value = 751.875
$ws.Range("I5").Value = 108.666664
$ws.Range("J5").Value = 1137.8
$ws.Range("K5").Value = 108.666664
$ws.Range("L5").Value = 1137.8
$ws.Range("M5").Value = 6.333336000000003
$ws.Range("N5").Value = -1367.8
$ws.Range("H6").Value = 474.5
$ws.Range("J6").Value = 349.5
$ws.Range("L6").Value = 1048.5
$ws.Range("N6").Value = -1272.5
$ws.Range("H18").Value = 7633
$ws.Range("J18").Value = 999
$ws.Range("L18").Value = 999
$ws.Range("N18").Value = -1567
$ws.Range("H40").Value = 1772.7693
$ws.Range("I40").Value = 1506.5333
$ws.Range("J40").Value = 2135.818
$ws.Range("K40").Value = 1506.5333
$ws.Range("L40").Value = 2135.818
$ws.Range("M40").Value = -1331.5333
$ws.Range("N40").Value = -2485.818
$ws.Range("H62").Value = 1744.5714
$ws.Range("I62").Value = 1732
$ws.Range("K62").Value = 1732
$ws.Range("M62").Value = -1108
$ws.Range("H65").Value = 1744.5714
$ws.Range("I65").Value = 1732
$ws.Range("K65").Value = 8660
$ws.Range("M65").Value = -5540
$ws.Range("H74").Value = 107778.47
$ws.Range("I74").Value = 129722.75
$ws.Range("K74").Value = 129722.75
$ws.Range("M74").Value = -128786.75
$ws.Range("H77").Value = 107778.47
$ws.Range("I77").Value = 129722.75
$ws.Range("K77").Value = 648613.75
$ws.Range("M77").Value = -643933.75
$ws.Range("H92").Value = 2107.8
$ws.Range("I92").Value = 3400
$ws.Range("K92").Value = 3400
$ws.Range("M92").Value = -2152
$ws.Range("H103").Value = 916.3333
$ws.Range("J103").Value = 1099.5
$ws.Range("L103").Value = 3298.5
$ws.Range("N103").Value = -4470.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 9849.75
$ws.Range("J125").Value = 18099.75
$ws.Range("L125").Value = 162897.75
$ws.Range("N125").Value = -167817.75
$ws.Range("H135").Value = 396.7143
$ws.Range("I135").Value = 396.7143
$ws.Range("K135").Value = 3570.4287
$ws.Range("M135").Value = -1035.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 1998.5
$ws.Range("I137").Value = 1998
$ws.Range("J137").Value = 1999
$ws.Range("K137").Value = 5994
$ws.Range("L137").Value = 5997
$ws.Range("M137").Value = -3444
$ws.Range("N137").Value = -11097
$ws.Range("H5").Value = 425.42856
$ws.Range("I5").Value = 462.66666
$ws.Range("J5").Value = 202
$ws.Range("K5").Value = 462.66666
$ws.Range("L5").Value = 202
$ws.Range("M5").Value = -350.66666
$ws.Range("N5").Value = -426
$ws.Range("H60").Value = 20000
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21466
$ws.Range("H63").Value = 3461.647
$ws.Range("I63").Value = 3553
$ws.Range("K63").Value = 3553
$ws.Range("M63").Value = -2867
$ws.Range("H66").Value = 3461.647
$ws.Range("I66").Value = 3553
$ws.Range("K66").Value = 17765
$ws.Range("M66").Value = -14333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 425.42856
$ws.Range("I4").Value = 462.66666
$ws.Range("J4").Value = 202
$ws.Range("K4").Value = 462.66666
$ws.Range("L4").Value = 202
$ws.Range("M4").Value = -347.66666
$ws.Range("N4").Value = -432
$ws.Range("H75").Value = 28129.166
$ws.Range("I75").Value = 28129.166
$ws.Range("K75").Value = 28129.166
$ws.Range("M75").Value = -27193.166
$ws.Range("H78").Value = 28129.166
$ws.Range("I78").Value = 28129.166
$ws.Range("K78").Value = 28129.166
$ws.Range("M78").Value = -79707.49800000001
$ws.Range("H86").Value = 1365.3334
$ws.Range("I86").Value = 1305.2
$ws.Range("K86").Value = 1305.2
$ws.Range("M86").Value = -182.2
$ws.Range("H89").Value = 1365.3334
$ws.Range("I89").Value = 1305.2
$ws.Range("K89").Value = 6526
$ws.Range("M89").Value = -910
$ws.Range("H94").Value = 644
$ws.Range("I94").Value = 592.75
$ws.Range("J94").Value = 849
$ws.Range("K94").Value = 592.75
$ws.Range("L94").Value = 849
$ws.Range("M94").Value = -141.75
$ws.Range("N94").Value = -1751
$ws.Range("H97").Value = 18487.2
$ws.Range("I97").Value = 18134
$ws.Range("K97").Value = 18134
$ws.Range("M97").Value = -17143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2574.4
$ws.Range("I31").Value = 799
$ws.Range("J31").Value = 3758
$ws.Range("K31").Value = 799
$ws.Range("L31").Value = 3758
$ws.Range("M31").Value = -504
$ws.Range("N31").Value = -4348
$ws.Range("H34").Value = 2574.4
$ws.Range("I34").Value = 799
$ws.Range("J34").Value = 3758
$ws.Range("K34").Value = 799
$ws.Range("L34").Value = 3758
$ws.Range("M34").Value = -597
$ws.Range("N34").Value = -4162
$ws.Range("H62").Value = 3900
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 3900
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H74").Value = 59499.5
$ws.Range("J74").Value = 59499.5
$ws.Range("L74").Value = 59499.5
$ws.Range("N74").Value = -61247.5
$ws.Range("H77").Value = 59499.5
$ws.Range("J77").Value = 59499.5
$ws.Range("L77").Value = 178498.5
$ws.Range("N77").Value = -187234.5
$ws.Range("H94").Value = 67965.94
$ws.Range("J94").Value = 3839.5
$ws.Range("L94").Value = 3839.5
$ws.Range("N94").Value = -4741.5
$ws.Range("H95").Value = 13945.75
$ws.Range("J95").Value = 13945.75
$ws.Range("L95").Value = 13945.75
$ws.Range("N95").Value = -19437.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 649.375
$ws.Range("J92").Value = 732.5
$ws.Range("L92").Value = 2197.5
$ws.Range("N92").Value = -4693.5
$ws.Range("H97").Value = 897.1111
$ws.Range("I97").Value = 1347.4
$ws.Range("K97").Value = 4042.2
$ws.Range("M97").Value = -3546.2
$ws.Range("H107").Value = 574.5
$ws.Range("I107").Value = 599.3333
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1797.9999
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 122.0001
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5304
$ws.Range("I70").Value = 5103.5
$ws.Range("K70").Value = 5103.5
$ws.Range("M70").Value = -4833.5
$ws.Range("H73").Value = 5304
$ws.Range("I73").Value = 5103.5
$ws.Range("K73").Value = 5103.5
$ws.Range("M73").Value = -4167.5
$ws.Range("H107").Value = 3088.3333
$ws.Range("I107").Value = 1447.5
$ws.Range("J107").Value = 3557.1428
$ws.Range("K107").Value = 1447.5
$ws.Range("L107").Value = 3557.1428
$ws.Range("M107").Value = 472.5
$ws.Range("N107").Value = -7397.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1190
$ws.Range("I9").Value = 1190
$ws.Range("K9").Value = 1190
$ws.Range("M9").Value = -966

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 22258
$ws.Range("J74").Value = 22208.5
$ws.Range("L74").Value = 22208.5
$ws.Range("N74").Value = -24080.5
$ws.Range("H77").Value = 22258
$ws.Range("J77").Value = 22208.5
$ws.Range("L77").Value = 66625.5
$ws.Range("N77").Value = -75985.5
$ws.Range("H81").Value = 1252677.8
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 1431417.4
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 2862834.8
$ws.Range("M81").Value = -1939
$ws.Range("N81").Value = -2864956.8
$ws.Range("H84").Value = 1252677.8
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 1431417.4
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 14314174
$ws.Range("M84").Value = -9696
$ws.Range("N84").Value = -14324782
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("I122").Value = 1466.5
$ws.Range("J122").Value = 1898.6
$ws.Range("K122").Value = 4399.5
$ws.Range("L122").Value = 5695.799999999999
$ws.Range("M122").Value = -1949.5
$ws.Range("N122").Value = -10595.8
$ws.Range("H132").Value = 5167.6665
$ws.Range("I132").Value = 5167.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15502.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12972.9995
$ws.Range("N132").ClearContents()
